$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.819.15"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.30%  '
$ws.Range('D3').Value = "'1.891.65"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'0.7782"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.86%  '
$ws.Range('D6').Value = "'243.51"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D8').Value = "'0.3125"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.93%  '
$ws.Range('D9').Value = "'25.25"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -5.76%  '
$ws.Range('D10').Value = "'0.07145"
$ws.Range('D10').ClearFormats()
$ws.Range('D11').Value = "'0.08079"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = "'0.7604"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.46%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = "'1.916.78"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'5.450"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.80%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = "'92.08"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.40%  '
$ws.Range('D16').Value = "'6.133"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.45%  '
$ws.Range('D17').Value = "'29.850.14"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').Value = "'13.91"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.40%  '
$ws.Range('D19').Value = "'242.99"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.97%  '
$ws.Range('D20').Value = "'0.000007763"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').Value = "'1.001"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').Value = "'2.125.03"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.44%  '
$ws.Range('D23').Value = "'8.066"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +11.60%  '
$ws.Range('D24').Value = "'1.001"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = "'0.1621"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.37%  '
$ws.Range('D26').Value = "'9.382"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').Value = "'162.12"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.31%  '
$ws.Range('E28').Value = '  -1.71%  '
$ws.Range('E29').Value = '  -3.61%  '
$ws.Range('E30').Value = '  +2.74%  '
$ws.Range('D31').Value = "'1.546"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('D32').Value = "'4.468"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('D33').Value = "'4.098"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('D34').Value = "'0.05520"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.15%  '
$ws.Range('D35').Value = "'1.262"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('D36').Value = "'0.7411"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').Value = "'0.9926"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.08%  '
$ws.Range('D38').Value = "'2.617"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.28%  '
$ws.Range('D39').Value = "'0.01910"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  -0.62%  '
$ws.Range('D41').Value = "'1.140.49"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +12.32%  '
$ws.Range('D42').Value = "'73.61"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').Value = "'0.4402"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('D44').Value = "'5.836"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.21%  '
$ws.Range('D45').Value = "'0.8500"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = "'1.001"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = "'103.51"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('D48').Value = "'1.865"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.59%  '
$ws.Range('D49').Value = "'9.913"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = "'7.427"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.29%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = "'3.015"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +10.67%  '
